# Improvements_to_MKPedals.xlsx edits
# - Mark several "Done?" checkboxes as TRUE
# - Replace the stray "Cre" row content (row 23) with real item-filter-view notes
# - Delete the final blank row (row 50)
# - Update the active selection in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tick the "Done?" checkboxes for the relevant improvements.
$ws.Range("D4").Value = $true
$ws.Range("D5").Value = $true
$ws.Range("D13").Value = $true
$ws.Range("D14").Value = $true
$ws.Range("D15").Value = $true
$ws.Range("D18").Value = $true

# Row 23 used to hold a stray, unfinished "Cre" entry. Replace it with the
# real "item filter view" improvement entry (Page / Improvement / Notes).
$ws.Range("A23").Value = "read.php/admin_read.php"
$ws.Range("B23").Value = "Improve the filter dropdown to dynamically update when new effects types are added"
$ws.Range("C23").Value = 'Could probably take all of the effects types under "effect_type" in the database, push them to an array and use a loop to echo them all to the dropdown list'

# That note text wraps onto several lines, so grow the row to fit it (matches
# the height Excel settles on for the other multi-line note rows, e.g. row 20).
$ws.Rows(23).RowHeight = 72.5

# Remove the trailing blank row at the bottom of the table.
$ws.Rows(50).Delete()

# Restore the cursor/scroll position left by the editing session.
$ws.Range("D13").Select() | Out-Null
